$d = $word.ActiveDocument

# The three new paragraphs must be inserted right before the final
# (bookmark-only) paragraph of the document -- i.e. right after the
# paragraph ending in "...saved in db". We build a zero-length Range
# at the exact character offset between those two paragraphs via
# $d.Range(pos, pos) -- this (unlike Duplicate+Collapse on the
# bookmark paragraph's own Range) leaves the bookmark markers and the
# preceding paragraph's runs untouched.
$paraCount = $d.Paragraphs.Count
$prevPara = $d.Paragraphs.Item($paraCount - 1)
$pos = $prevPara.Range.End
$r = $d.Range($pos, $pos)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="HTMLPreformatted"/>
    <w:spacing w:line="276" w:lineRule="auto"/>
    <w:ind w:left="720"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="HTMLPreformatted"/>
    <w:spacing w:line="276" w:lineRule="auto"/>
    <w:ind w:left="720"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
    <w:t>FRONTEND</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="HTMLPreformatted"/>
    <w:spacing w:line="276" w:lineRule="auto"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body></w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)

$d.Save()
